# Generate Report for Handback
# Update the generated/handoff/handback timestamps shown on the report sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 09:34:56"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 09:34:44"
$wsZhCn.Range("K2").Value = "2016-09-07 09:35:46"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-07 09:34:56"
$wsDeDe.Range("K2").Value = "2016-09-07 09:36:14"
